# Normalize the "Recorded By" (column G) attribution lists so that the
# synthetic "System" actor is always listed first among the recorders,
# matching the canonical ordering used by the sync job. A couple of rows
# (no "System" entry) instead get their two human recorders swapped to
# the same canonical order (dnasr281@gmail.com before admin@admin.com).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

$colG = 7

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $colG)
    $raw = $cell.Value2

    if ($null -eq $raw) { continue }
    $current = "$raw"
    if ($current -eq "") { continue }

    $parts = @($current.Split(",") | ForEach-Object { $_.Trim() })

    $newValue = $null

    if ($parts -contains "System") {
        $rest = @($parts | Where-Object { $_ -ne "System" })
        $reordered = @("System") + $rest
        $newValue = [string]::Join(", ", $reordered)
    }
    elseif ($current -eq "admin@admin.com, dnasr281@gmail.com") {
        $newValue = "dnasr281@gmail.com, admin@admin.com"
    }

    if ($null -ne $newValue -and $newValue -ne $current) {
        $cell.Value = $newValue
    }
}
